$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "2025-04-28 12:24:12"
$ws.Range("B14").Value = 221
